$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header texts (B1, C1) - shared string updates
$ws.Range("B1").Value = "AVERAGE_TIME_PER_ISSUES_WO_FT"
$ws.Range("C1").Value = "AVERAGE_TIME_PER_ISSUES_WITH_FT"

# Update the selection shown in the sheet view to B1:C1
$ws.Range("B1:C1").Select()
